$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-9 from 45207 to 45208
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 3).Value = 45208
}
